$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.088.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.302.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.509"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.501"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.32%  "
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.04%  "
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.45%  "
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.658.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.299.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.811"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.014.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.86%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.29%  "
$ws.Range("E31").Value = "  -5.79%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.70"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.52%  "
$ws.Range("E35").Value = "  -1.15%  "
$ws.Range("E36").Value = "  -1.11%  "
$ws.Range("E37").Value = "  +2.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0695"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("E41").Value = "  -1.61%  "
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.982.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.25%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("E48").Value = "  -0.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.526.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.74%  "
